$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Delete obsolete worker rows (SILVIA row16, LULY rows 17-18) - remaining ANA rows shift up
$ws.Rows("16:18").Delete()

# 2. Swap period/values between (now) row16 and row18 so periods run ascending 2204,2205,2206
$e16 = $ws.Range("E16").Value()
$f16 = $ws.Range("F16").Value()
$e18 = $ws.Range("E18").Value()
$f18 = $ws.Range("F18").Value()

$ws.Range("E16").Value = $e18
$ws.Range("F16").Value = $f18
$ws.Range("E18").Value = $e16
$ws.Range("F18").Value = $f16

# 3. Update the summary fields: total overdue value, worker count, period count
$ws.Range("E11").Value = 159120
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 3

Write-Host "done rows"
